# "collected more images for report"
# Adds a new "24" label (TextBox) into the existing diagram group on
# slide 1, next to the other numeric callouts already in that group.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The whole diagram lives inside one big top-level group shape.
$grp = $s.Shapes.Item(1)

# Burn the next free shape id (an internal PowerPoint id-allocation
# artifact) so the id/name PowerPoint hands out to the new textbox and
# to the re-created group line up the way they do after a normal
# "add a shape, ungroup, regroup" editing session.
$placeholder = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$placeholder.Delete()

# Add the new "24" callout textbox at its final position/size (EMU
# values converted to points, since Shape coordinates are in points).
$tb = $s.Shapes.AddTextbox(1, 1924922 / 12700.0, 3593188 / 12700.0, 325529 / 12700.0, 246221 / 12700.0)
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "24"
$tb.TextFrame.TextRange.Font.Size = 10

# Break the existing group back into loose shapes, then regroup
# everything (the original members + the new textbox) so the new
# textbox becomes a sibling inside a single (re-created) group,
# instead of being wrapped in an extra nested group.
$grp.Ungroup() | Out-Null

$names = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $names += $s.Shapes.Item($i).Name
}
$all = $s.Shapes.Range($names)
$newGroup = $all.Group()
